$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (name unchanged)
$ws.Range("B2").Value = 0.4213705682003807
$ws.Range("C2").Value = 0.4213705682003805
$ws.Range("D2").Value = 0.4213705682003804

# Row 3 - RandomForestRegressor (name unchanged)
$ws.Range("B3").Value = 0.0324731432071553
$ws.Range("C3").Value = 0.03242383074417194
$ws.Range("D3").Value = 0.1699595858781286

# Row 4 - name changes from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.03513881945894168
$ws.Range("C4").Value = 0.03459865781640419
$ws.Range("D4").Value = 0.1347765933567942

# Row 5 - name changes from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.03073366596206964
$ws.Range("C5").Value = 0.02938574997184816
$ws.Range("D5").Value = 0.05941973989784682
